$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Common")

# Insert a new row before row 70 (shifts rows 70-129 down to 71-130)
$ws.Rows.Item(70).Insert()

# Set the new row's content
$ws.Range("A70").Value = "VSTAT License File"

# Add the comment describing the new license file field
$ws.Range("A70").AddComment("Optional License File for Elasticsearch [default: ]")
